$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97; this shifts the existing rows 97-155
# down to 98-156 (and Excel auto-extends the used range / dimension).
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new "Haba" price record.
$ws.Cells.Item(97, 1).Value = 9
$ws.Cells.Item(97, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(97, 3).Value = "Metropolitana"
$ws.Cells.Item(97, 4).Value = 44488
$ws.Cells.Item(97, 5).Value = 13
$ws.Cells.Item(97, 6).Value = 100112026
$ws.Cells.Item(97, 7).Value = "Haba"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 52
$ws.Cells.Item(97, 11).Value = 6000
$ws.Cells.Item(97, 12).Value = 7000
$ws.Cells.Item(97, 13).Value = 6500
$ws.Cells.Item(97, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(97, 16).Value = 260
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date-formatted style used by the
# rest of column D (style index 2 in the original workbook).
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
